# Update manuscript model results table with refreshed statistics.
# Purely-numeric-looking values are written with a leading apostrophe
# (forcing text) and the cell style is reset back to "Normal" afterwards
# so the cells stay plain shared-string text cells (matching the source
# table where every value, even number-like ones, is stored as text).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Intercept row (row 2): coef.final / CI.full / coef.final / CI.final
$ws.Range("B2").Value = "'2.52"
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").Value = "[1.82; 3.48]"
$ws.Range("F2").Value = "'1.86"
$ws.Range("F2").Style = "Normal"
$ws.Range("G2").Value = "[1.36; 2.54]"

# dPTH row (row 3): CI.full, coef.simple
$ws.Range("C3").Value = "[1.05; 1.12]"
$ws.Range("J3").Value = "'1.09"
$ws.Range("J3").Style = "Normal"

# CorrCa24u row (row 4): CI.full, CI.final
$ws.Range("C4").Value = "[1.13; 1.92]"
$ws.Range("G4").Value = "[1.11; 1.87]"

# BSKgezien row (row 5): coef.full, CI.full, coef.final, CI.final
$ws.Range("B5").Value = "'3.67"
$ws.Range("B5").Style = "Normal"
$ws.Range("C5").Value = "[1.46; 9.23]"
$ws.Range("F5").Value = "'3.94"
$ws.Range("F5").Style = "Normal"
$ws.Range("G5").Value = "[1.63; 9.53]"

# Age_Years row (row 6): CI.full
$ws.Range("C6").Value = "[0.61; 2.11]"

# Sex row (row 7): CI.full
$ws.Range("C7").Value = "[0.45; 2.02]"

# surgery_type row (row 8): CI.full
$ws.Range("C8").Value = "[0.44; 4.91]"

# CHKD row (row 9): CI.full
$ws.Range("C9").Value = "[0.61; 2.87]"

# C-index row (row 10): CI.simple
$ws.Range("K10").Value = "[0.81; 0.90]"
